$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new header cells (G1, H1)
$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"

# Match formatting of the existing header row (bold, bordered, centered)
$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)  # xlPasteFormats

# Add the new data values (G2, H2)
$ws.Range("G2").Value = 0.1258822953001072
$ws.Range("H2").Value = 0.988
